$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Remove the second row entirely (value + style), shrinking the used range to A1:B1
$ws.Range("A2:B2").Clear()

# Overwrite row 1 with new short values. Write B1 before A1 so the shared-string
# table ends up ordered [0]="B", [1]="A" (matching the target sharedStrings.xml)
$ws.Range("B1").Value = "B"
$ws.Range("A1").Value = "A"

# Columns no longer need to "best fit" the old long header text; give both
# columns a plain custom width
$ws.Range("A1:B1").ColumnWidth = 8.6

# Move the selection to C2
$ws.Range("C2").Select()

# Restore the originally active sheet/tab (selecting on "Global" above activated it)
$wb.Worksheets.Item("Action1").Activate()
